$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

function Set-TextValue($Range, $Text) {
    # Force the cell to be stored as text (matching the existing shared-string
    # cell type) instead of letting Excel auto-convert the numeric-looking
    # string into a Number, then strip the temporary Text number format back
    # off so the cell's style is left exactly as it was before the edit.
    $Range.NumberFormat = "@"
    $Range.Value = $Text
    $Range.ClearFormats()
}

# Enterprises density (per 1000 people)
Set-TextValue -Range $ws.Range("B13") -Text "20.72"
Set-TextValue -Range $ws.Range("C13") -Text "1.12"
Set-TextValue -Range $ws.Range("D13") -Text "21.84"

# Enterprises (% of total)
Set-TextValue -Range $ws.Range("B16") -Text "94.79"
Set-TextValue -Range $ws.Range("C16") -Text "5.14"
Set-TextValue -Range $ws.Range("D16") -Text "99.92"

Write-Output "Updated Tajikistan enterprise density/percentage figures"
